# Auto-generated Excel COM-interop script applying the Mateus_Profits.xlsx value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets (crafting profit recalculation).
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 260.75
$ws.Range("I38").Value = 260.75
$ws.Range("K38").Value = 782.25
$ws.Range("M38").Value = -410.25
$ws.Range("H74").Value = 9488.739
$ws.Range("I74").Value = 7314.25
$ws.Range("J74").Value = 11860.909
$ws.Range("K74").Value = 7314.25
$ws.Range("L74").Value = 11860.909
$ws.Range("M74").Value = -6378.25
$ws.Range("N74").Value = -13732.909
$ws.Range("H77").Value = 9488.739
$ws.Range("I77").Value = 7314.25
$ws.Range("J77").Value = 11860.909
$ws.Range("K77").Value = 36571.25
$ws.Range("L77").Value = 59304.545
$ws.Range("M77").Value = -31891.25
$ws.Range("N77").Value = -68664.545
$ws.Range("H98").Value = 1431.5
$ws.Range("J98").Value = 1990
$ws.Range("L98").Value = 1990
$ws.Range("N98").Value = -4986
$ws.Range("H111").Value = 398.5
$ws.Range("I111").Value = 249.5
$ws.Range("K111").Value = 748.5
$ws.Range("M111").Value = 2318.5
$ws.Range("H116").Value = 3999.5454
$ws.Range("I116").Value = 3227.8572
$ws.Range("K116").Value = 3227.8572
$ws.Range("M116").Value = 214.1428000000001
$ws.Range("H122").Value = 1431.5
$ws.Range("J122").Value = 1990
$ws.Range("L122").Value = 5970
$ws.Range("N122").Value = -10870
$ws.Range("H131").Value = 3861.3125
$ws.Range("I131").Value = 2785.4
$ws.Range("K131").Value = 8356.200000000001
$ws.Range("M131").Value = -3316.200000000001
$ws.Range("H132").Value = 2117.2173
$ws.Range("I132").Value = 1973.4736
$ws.Range("K132").Value = 5920.4208
$ws.Range("M132").Value = -3390.4208
$ws.Range("H137").Value = 1138.2727
$ws.Range("I137").Value = 1080.0938
$ws.Range("K137").Value = 3240.2814
$ws.Range("M137").Value = -690.2814000000003

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4962.909
$ws.Range("I2").Value = 3699.25
$ws.Range("K2").Value = 3699.25
$ws.Range("M2").Value = -3586.25
$ws.Range("H32").Value = 12634.772
$ws.Range("I32").Value = 10408.941
$ws.Range("K32").Value = 10408.941
$ws.Range("M32").Value = -10121.941
$ws.Range("H74").Value = 3094.2122
$ws.Range("I74").Value = 2326.0967
$ws.Range("K74").Value = 2326.0967
$ws.Range("M74").Value = -1452.0967
$ws.Range("H77").Value = 3094.2122
$ws.Range("I77").Value = 2326.0967
$ws.Range("K77").Value = 11630.4835
$ws.Range("M77").Value = -7262.4835
$ws.Range("H88").Value = 2407.6155
$ws.Range("I88").Value = 2375.375
$ws.Range("J88").Value = 2459.2
$ws.Range("K88").Value = 2375.375
$ws.Range("L88").Value = 2459.2
$ws.Range("M88").Value = -1969.375
$ws.Range("N88").Value = -3271.2
$ws.Range("H91").Value = 2407.6155
$ws.Range("I91").Value = 2375.375
$ws.Range("J91").Value = 2459.2
$ws.Range("K91").Value = 2375.375
$ws.Range("L91").Value = 2459.2
$ws.Range("M91").Value = -971.375
$ws.Range("N91").Value = -5267.2
$ws.Range("H102").Value = 3016.84
$ws.Range("I102").Value = 2214.8
$ws.Range("K102").Value = 2214.8
$ws.Range("M102").Value = -592.8000000000002
$ws.Range("H116").Value = 4962.909
$ws.Range("I116").Value = 3699.25
$ws.Range("K116").Value = 3699.25
$ws.Range("M116").Value = -1405.25
$ws.Range("H122").Value = 2994.3333
$ws.Range("I122").Value = 2994.3333
$ws.Range("K122").Value = 8982.999899999999
$ws.Range("M122").Value = -6532.999899999999
$ws.Range("H132").Value = 1535.6123
$ws.Range("I132").Value = 1550.9584
$ws.Range("J132").Value = 799
$ws.Range("K132").Value = 4652.8752
$ws.Range("L132").Value = 2397
$ws.Range("M132").Value = -2122.8752
$ws.Range("N132").Value = -7457
$ws.Range("H135").Value = 63574.816
$ws.Range("J135").Value = 63574.816
$ws.Range("L135").Value = 63574.816
$ws.Range("N135").Value = -73714.81599999999

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4962.909
$ws.Range("I3").Value = 3699.25
$ws.Range("K3").Value = 3699.25
$ws.Range("M3").Value = -3585.25
$ws.Range("H20").Value = 2811.3928
$ws.Range("J20").Value = 3066.1538
$ws.Range("L20").Value = 3066.1538
$ws.Range("N20").Value = -3560.1538
$ws.Range("H86").Value = 1434.1923
$ws.Range("I86").Value = 1249.25
$ws.Range("J86").Value = 2050.6667
$ws.Range("K86").Value = 1249.25
$ws.Range("L86").Value = 2050.6667
$ws.Range("M86").Value = -126.25
$ws.Range("N86").Value = -4296.6667
$ws.Range("H89").Value = 1434.1923
$ws.Range("I89").Value = 1249.25
$ws.Range("J89").Value = 2050.6667
$ws.Range("K89").Value = 6246.25
$ws.Range("L89").Value = 10253.3335
$ws.Range("M89").Value = -630.25
$ws.Range("N89").Value = -21485.3335
$ws.Range("H94").Value = 2224.75
$ws.Range("I94").Value = 1632.6666
$ws.Range("J94").Value = 2580
$ws.Range("K94").Value = 1632.6666
$ws.Range("L94").Value = 2580
$ws.Range("M94").Value = -1181.6666
$ws.Range("N94").Value = -3482
$ws.Range("H132").Value = 97749
$ws.Range("J132").Value = 97749
$ws.Range("L132").Value = 97749
$ws.Range("N132").Value = -107869
$ws.Range("H140").Value = 94666.336
$ws.Range("J140").Value = 94666.336
$ws.Range("L140").Value = 94666.336
$ws.Range("N140").Value = -105026.336

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 46196.285
$ws.Range("J112").Value = 46196.285
$ws.Range("L112").Value = 46196.285
$ws.Range("N112").Value = -49150.285
$ws.Range("H122").Value = 3027.4614
$ws.Range("I122").Value = 2987.647
$ws.Range("K122").Value = 8962.940999999999
$ws.Range("M122").Value = -6512.940999999999
$ws.Range("H134").Value = 5166.6113
$ws.Range("I134").Value = 3033.25
$ws.Range("K134").Value = 9099.75
$ws.Range("M134").Value = -6564.75
$ws.Range("H141").Value = 178392
$ws.Range("J141").Value = 178392
$ws.Range("L141").Value = 178392
$ws.Range("N141").Value = -188752

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1793.5
$ws.Range("I13").Value = 40
$ws.Range("J13").Value = 2144.2
$ws.Range("K13").Value = 120
$ws.Range("L13").Value = 6432.599999999999
$ws.Range("M13").Value = 48
$ws.Range("N13").Value = -6768.599999999999
$ws.Range("H64").Value = 1763.9
$ws.Range("I64").Value = 1626.5555
$ws.Range("K64").Value = 4879.666499999999
$ws.Range("M64").Value = -4609.666499999999
$ws.Range("H67").Value = 1763.9
$ws.Range("I67").Value = 1626.5555
$ws.Range("K67").Value = 4879.666499999999
$ws.Range("M67").Value = -3943.666499999999
$ws.Range("H88").Value = 14962.25
$ws.Range("J88").Value = 14962.25
$ws.Range("L88").Value = 44886.75
$ws.Range("N88").Value = -45742.75
$ws.Range("H91").Value = 14962.25
$ws.Range("J91").Value = 14962.25
$ws.Range("L91").Value = 44886.75
$ws.Range("N91").Value = -47850.75
$ws.Range("H94").Value = 2740.6667
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 3611
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 10833
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -12185
$ws.Range("H113").Value = 500
$ws.Range("J113").Value = 500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840
$ws.Range("H121").Value = 167800.33
$ws.Range("I121").Value = 442.25
$ws.Range("J121").Value = 502516.5
$ws.Range("K121").Value = 1326.75
$ws.Range("L121").Value = 1507549.5
$ws.Range("M121").Value = -16.75
$ws.Range("N121").Value = -1510169.5
$ws.Range("H129").Value = 557356.5600000001
$ws.Range("I129").Value = 143910.28
$ws.Range("J129").Value = 779981.4399999999
$ws.Range("K129").Value = 431730.84
$ws.Range("L129").Value = 2339944.32
$ws.Range("M129").Value = -426730.84
$ws.Range("N129").Value = -2349944.32

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15166.667
$ws.Range("J70").Value = 16200
$ws.Range("L70").Value = 16200
$ws.Range("N70").Value = -16740
$ws.Range("H73").Value = 15166.667
$ws.Range("J73").Value = 16200
$ws.Range("L73").Value = 16200
$ws.Range("N73").Value = -18072
$ws.Range("H80").Value = 2691.7778
$ws.Range("I80").Value = 1921.7142
$ws.Range("J80").Value = 3181.818
$ws.Range("K80").Value = 1921.7142
$ws.Range("L80").Value = 3181.818
$ws.Range("M80").Value = -923.7141999999999
$ws.Range("N80").Value = -5177.818
$ws.Range("H83").Value = 2691.7778
$ws.Range("I83").Value = 1921.7142
$ws.Range("J83").Value = 3181.818
$ws.Range("K83").Value = 9608.571
$ws.Range("L83").Value = 15909.09
$ws.Range("M83").Value = -4616.571
$ws.Range("N83").Value = -25893.09
$ws.Range("H122").Value = 2341.2222
$ws.Range("I122").Value = 2346.8333
$ws.Range("J122").Value = 2330
$ws.Range("K122").Value = 7040.499899999999
$ws.Range("L122").Value = 6990
$ws.Range("M122").Value = -4590.499899999999
$ws.Range("N122").Value = -11890
$ws.Range("H132").Value = 3195.9268
$ws.Range("I132").Value = 2846.359
$ws.Range("K132").Value = 8539.076999999999
$ws.Range("M132").Value = -6009.076999999999

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5612.1113
$ws.Range("I7").Value = 5480.375
$ws.Range("K7").Value = 5480.375
$ws.Range("M7").Value = -5368.375
$ws.Range("H40").Value = 4000.9565
$ws.Range("I40").Value = 3906.1538
$ws.Range("J40").Value = 4124.2
$ws.Range("K40").Value = 3906.1538
$ws.Range("L40").Value = 4124.2
$ws.Range("M40").Value = -3770.1538
$ws.Range("N40").Value = -4396.2
$ws.Range("H46").Value = 2898.3333
$ws.Range("I46").Value = 2898.3333
$ws.Range("K46").Value = 2898.3333
$ws.Range("M46").Value = -2710.3333
$ws.Range("H55").Value = 35714390
$ws.Range("I55").Value = 50000104
$ws.Range("J55").Value = 116.25
$ws.Range("K55").Value = 50000104
$ws.Range("L55").Value = 116.25
$ws.Range("M55").Value = -49999931
$ws.Range("N55").Value = -462.25
$ws.Range("H61").Value = 30199.914
$ws.Range("I61").Value = 45433.61
$ws.Range("J61").Value = 1002
$ws.Range("K61").Value = 45433.61
$ws.Range("L61").Value = 1002
$ws.Range("M61").Value = -45231.61
$ws.Range("N61").Value = -1406
$ws.Range("H68").Value = 5645
$ws.Range("I68").Value = 2241.8333
$ws.Range("J68").Value = 10749.75
$ws.Range("K68").Value = 2241.8333
$ws.Range("L68").Value = 10749.75
$ws.Range("M68").Value = -1492.8333
$ws.Range("N68").Value = -12247.75
$ws.Range("H71").Value = 5645
$ws.Range("I71").Value = 2241.8333
$ws.Range("J71").Value = 10749.75
$ws.Range("K71").Value = 11209.1665
$ws.Range("L71").Value = 53748.75
$ws.Range("M71").Value = -7465.166499999999
$ws.Range("N71").Value = -61236.75
$ws.Range("H82").Value = 1240.2307
$ws.Range("I82").Value = 1146
$ws.Range("J82").Value = 1321
$ws.Range("K82").Value = 1146
$ws.Range("L82").Value = 1321
$ws.Range("M82").Value = -785
$ws.Range("N82").Value = -2043
$ws.Range("H85").Value = 1240.2307
$ws.Range("I85").Value = 1146
$ws.Range("J85").Value = 1321
$ws.Range("K85").Value = 1146
$ws.Range("L85").Value = 1321
$ws.Range("M85").Value = 102
$ws.Range("N85").Value = -3817
$ws.Range("H93").Value = 16083.294
$ws.Range("I93").Value = 6192.6
$ws.Range("J93").Value = 30212.857
$ws.Range("K93").Value = 6192.6
$ws.Range("L93").Value = 30212.857
$ws.Range("M93").Value = -4944.6
$ws.Range("N93").Value = -32708.857
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents() | Out-Null
$ws.Range("H113").Value = 30199.914
$ws.Range("I113").Value = 45433.61
$ws.Range("J113").Value = 1002
$ws.Range("K113").Value = 45433.61
$ws.Range("L113").Value = 1002
$ws.Range("M113").Value = -43263.61
$ws.Range("N113").Value = -5342
$ws.Range("H126").Value = 5612.1113
$ws.Range("I126").Value = 5480.375
$ws.Range("K126").Value = 16441.125
$ws.Range("M126").Value = -13971.125

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 84235.5
$ws.Range("J46").Value = 84235.5
$ws.Range("L46").Value = 84235.5
$ws.Range("N46").Value = -84697.5
$ws.Range("H62").Value = 4140.1
$ws.Range("I62").Value = 4225.25
$ws.Range("J62").Value = 4083.3333
$ws.Range("K62").Value = 4225.25
$ws.Range("L62").Value = 4083.3333
$ws.Range("M62").Value = -3601.25
$ws.Range("N62").Value = -5331.3333
$ws.Range("H65").Value = 4140.1
$ws.Range("I65").Value = 4225.25
$ws.Range("J65").Value = 4083.3333
$ws.Range("K65").Value = 21126.25
$ws.Range("L65").Value = 20416.6665
$ws.Range("M65").Value = -18006.25
$ws.Range("N65").Value = -26656.6665
$ws.Range("H81").Value = 1789.5
$ws.Range("I81").Value = 1888.3334
$ws.Range("J81").Value = 900
$ws.Range("K81").Value = 3776.6668
$ws.Range("L81").Value = 1800
$ws.Range("M81").Value = -2715.6668
$ws.Range("N81").Value = -3922
$ws.Range("H84").Value = 1789.5
$ws.Range("I84").Value = 1888.3334
$ws.Range("J84").Value = 900
$ws.Range("K84").Value = 18883.334
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = -13579.334
$ws.Range("N84").Value = -19608
$ws.Range("H100").Value = 878.25
$ws.Range("I100").Value = 743.2222
$ws.Range("K100").Value = 1486.4444
$ws.Range("M100").Value = -945.4444000000001
$ws.Range("H103").Value = 44901
$ws.Range("J103").Value = 44901
$ws.Range("L103").Value = 44901
$ws.Range("N103").Value = -47245
$ws.Range("H113").Value = 721.6
$ws.Range("I113").Value = 591.0625
$ws.Range("K113").Value = 1773.1875
$ws.Range("M113").Value = 396.8125
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H134").Value = 84235.5
$ws.Range("J134").Value = 84235.5
$ws.Range("L134").Value = 252706.5
$ws.Range("N134").Value = -257776.5
$ws.Range("H136").Value = 5557.579
$ws.Range("I136").Value = 5037.25
$ws.Range("K136").Value = 15111.75
$ws.Range("M136").Value = -12561.75
